# Update mappings.xlsx: insert a new "SupplyLookupMappings" sheet (prefix -> supply
# company lookup table) between "SupplyMappings" and "ProductMappings", and make it
# the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet right after "SupplyMappings".
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("SupplyMappings")
$lookup = $wb.Worksheets.Add($null, $afterSheet)
$lookup.Name = "SupplyLookupMappings"

# ---------------------------------------------------------------------------
# 2. Populate the prefix -> supply lookup table.
# ---------------------------------------------------------------------------
$data = @(
    @('prefix', 'supply'),
    @('BP-KM', 'BP'),
    @('PH-PH', 'Phillips 66'),
    @('FH-MG', 'Flint Hills'),
    @('GMK-MG', 'Growmark'),
    @('SC-HEP', 'Shell'),
    @('MPC', 'Marathon Petroleum'),
    @('JDS', 'JDS Energy'),
    @('SINCLAIR', 'Sinclair'),
    @('CVR', 'CVR Energy'),
    @('HTP', 'HTP Energy'),
    @('QUIKTRIP', 'QuikTrip'),
    @('SHELL', 'Shell'),
    @('WALLIS', 'Wallis Oil'),
    @('GROWMARK', 'Growmark'),
    @('MUSKET', 'Musket'),
    @('PH-KM', 'Phillips 66'),
    @('GMK-KM', 'Growmark'),
    @('UNBRANDED', 'Unbranded'),
    @('CHS', 'CHS'),
    @('P66', 'Phillips 66'),
    @('PSX', 'Phillips 66'),
    @('VALERO', 'Valero'),
    @('PH', 'Phillips 66'),
    @('PHI', 'Phillips 66'),
    @('PHX', 'Phillips 66'),
    @('PETRO-CAN', 'Petro-Canada'),
    @('CONOCO', 'Conoco'),
    @('CENEX', 'CHS'),
    @('BRT', 'Brentwood Oil'),
    @('GEMINI', 'Gemini Transport'),
    @('MIDCO', 'Midcoast Energy')
)

$r = 1
foreach ($pair in $data) {
    $lookup.Cells.Item($r, 1).Value2 = $pair[0]
    $lookup.Cells.Item($r, 2).Value2 = $pair[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3. Style the header row like the other lookup sheets: bold Cambria,
#    centered/top aligned, thin border.
# ---------------------------------------------------------------------------
$header = $lookup.Range("A1:B1")
$header.Font.Bold = $true
$header.Font.Name = "Cambria"
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 4. Column widths to roughly match the source layout.
# ---------------------------------------------------------------------------
$lookup.Columns.Item(1).ColumnWidth = 11.57
$lookup.Columns.Item(2).ColumnWidth = 18.38

# ---------------------------------------------------------------------------
# 5. Make the new sheet the active tab (matches activeTab="2" in the workbook).
# ---------------------------------------------------------------------------
$lookup.Activate()
$lookup.Range("D11").Select()

# ---------------------------------------------------------------------------
# 6. Scroll the previously-active "SupplierMappings" sheet so its stored
#    topLeftCell follows the source edit, while keeping its prior selection.
# ---------------------------------------------------------------------------
$supplier = $wb.Worksheets.Item("SupplierMappings")
$supplier.Activate()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 3
$supplier.Range("E11").Select()

# Leave the new sheet as the active / selected sheet, as in the target workbook.
$lookup.Activate()
